$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (pushes current rows 11-26 down to 12-27,
# inheriting the date-formatted style from the row that follows).
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly price record.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 45243
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112001
$ws.Range("G11").Value = "Berenjena"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 7500
$ws.Range("N11").Value = "`$/caja 60 unidades"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 125
$ws.Range("Q11").Value = 60
$ws.Range("R11").Value = "Hortaliza"
